# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) for the rows whose underlying Universalis price
# data changed, across all eight crafter sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1306.1072
$ws.Cells.Item(112, 9).Value = 473.33334
$ws.Cells.Item(112, 10).Value = 1406.04
$ws.Cells.Item(112, 11).Value = 1420.00002
$ws.Cells.Item(112, 12).Value = 4218.12
$ws.Cells.Item(112, 13).Value = -312.0000199999999
$ws.Cells.Item(112, 14).Value = -6434.12

$ws.Cells.Item(129, 8).Value = 851.13336
$ws.Cells.Item(129, 9).Value = 319.4
$ws.Cells.Item(129, 10).Value = 957.48
$ws.Cells.Item(129, 11).Value = 958.1999999999999
$ws.Cells.Item(129, 12).Value = 2872.44
$ws.Cells.Item(129, 13).Value = 4041.8
$ws.Cells.Item(129, 14).Value = -12872.44

$ws.Cells.Item(134, 8).Value = 51113.25
$ws.Cells.Item(134, 10).Value = 51113.25
$ws.Cells.Item(134, 12).Value = 51113.25
$ws.Cells.Item(134, 14).Value = -61253.25

$ws.Cells.Item(137, 8).Value = 1831.2115
$ws.Cells.Item(137, 9).Value = 1440.2812
$ws.Cells.Item(137, 10).Value = 2456.7
$ws.Cells.Item(137, 11).Value = 4320.8436
$ws.Cells.Item(137, 12).Value = 7370.099999999999
$ws.Cells.Item(137, 13).Value = -1770.8436
$ws.Cells.Item(137, 14).Value = -12470.1

$ws.Cells.Item(138, 8).Value = 2988.137
$ws.Cells.Item(138, 9).Value = 1355.6
$ws.Cells.Item(138, 10).Value = 4966.9697
$ws.Cells.Item(138, 11).Value = 4066.8
$ws.Cells.Item(138, 12).Value = 14900.9091
$ws.Cells.Item(138, 13).Value = 1073.2
$ws.Cells.Item(138, 14).Value = -25180.9091

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(94, 8).Value = 20000
$ws.Cells.Item(94, 10).Value = 20000
$ws.Cells.Item(94, 12).Value = 20000
$ws.Cells.Item(94, 14).Value = -21802

$ws.Cells.Item(110, 8).Value = 1696.8334
$ws.Cells.Item(110, 9).Value = 1518.8572
$ws.Cells.Item(110, 11).Value = 1518.8572
$ws.Cells.Item(110, 13).Value = 526.1428000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 16111.897
$ws.Cells.Item(134, 9).Value = 1161.2181
$ws.Cells.Item(134, 10).Value = 79364.766
$ws.Cells.Item(134, 11).Value = 3483.6543
$ws.Cells.Item(134, 12).Value = 238094.298
$ws.Cells.Item(134, 13).Value = -948.6543000000001
$ws.Cells.Item(134, 14).Value = -243164.298

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2094.8628
$ws.Cells.Item(31, 9).Value = 1465.8049
$ws.Cells.Item(31, 10).Value = 4674
$ws.Cells.Item(31, 11).Value = 1465.8049
$ws.Cells.Item(31, 12).Value = 4674
$ws.Cells.Item(31, 13).Value = -1170.8049
$ws.Cells.Item(31, 14).Value = -5264

$ws.Cells.Item(34, 8).Value = 2094.8628
$ws.Cells.Item(34, 9).Value = 1465.8049
$ws.Cells.Item(34, 10).Value = 4674
$ws.Cells.Item(34, 11).Value = 1465.8049
$ws.Cells.Item(34, 12).Value = 4674
$ws.Cells.Item(34, 13).Value = -1263.8049
$ws.Cells.Item(34, 14).Value = -5078

$ws.Cells.Item(58, 8).Value = 1058529.5
$ws.Cells.Item(58, 9).Value = 1378363.9
$ws.Cells.Item(58, 10).Value = 3076.3
$ws.Cells.Item(58, 11).Value = 1378363.9
$ws.Cells.Item(58, 12).Value = 3076.3
$ws.Cells.Item(58, 13).Value = -1378160.9
$ws.Cells.Item(58, 14).Value = -3482.3

$ws.Cells.Item(122, 8).Value = 8708.723
$ws.Cells.Item(122, 9).Value = 4839.1113
$ws.Cells.Item(122, 10).Value = 12578.333
$ws.Cells.Item(122, 11).Value = 14517.3339
$ws.Cells.Item(122, 12).Value = 37734.999
$ws.Cells.Item(122, 13).Value = -12067.3339
$ws.Cells.Item(122, 14).Value = -42634.999

$ws.Cells.Item(132, 8).Value = 2556.9348
$ws.Cells.Item(132, 9).Value = 2671.7314
$ws.Cells.Item(132, 10).Value = 2249.28
$ws.Cells.Item(132, 11).Value = 8015.1942
$ws.Cells.Item(132, 12).Value = 6747.84
$ws.Cells.Item(132, 13).Value = -5485.1942
$ws.Cells.Item(132, 14).Value = -11807.84

$ws.Cells.Item(136, 8).Value = 1058529.5
$ws.Cells.Item(136, 9).Value = 1378363.9
$ws.Cells.Item(136, 10).Value = 3076.3
$ws.Cells.Item(136, 11).Value = 4135091.7
$ws.Cells.Item(136, 12).Value = 9228.900000000001
$ws.Cells.Item(136, 13).Value = -4132541.7
$ws.Cells.Item(136, 14).Value = -14328.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 27783928
$ws.Cells.Item(5, 9).Value = 565.6667
$ws.Cells.Item(5, 10).Value = 55567292
$ws.Cells.Item(5, 11).Value = 1697.0001
$ws.Cells.Item(5, 12).Value = 166701876
$ws.Cells.Item(5, 13).Value = -1585.0001
$ws.Cells.Item(5, 14).Value = -166702100

$ws.Cells.Item(29, 8).Value = 600
$ws.Cells.Item(29, 10).Value = 600
$ws.Cells.Item(29, 12).Value = 1800
$ws.Cells.Item(29, 14).Value = -2354

$ws.Cells.Item(37, 8).Value = 76999.71000000001
$ws.Cells.Item(37, 10).Value = 76999.71000000001
$ws.Cells.Item(37, 12).Value = 230999.13
$ws.Cells.Item(37, 14).Value = -231223.13

$ws.Cells.Item(122, 8).Value = 640.3077
$ws.Cells.Item(122, 9).Value = 398.11765
$ws.Cells.Item(122, 11).Value = 3583.05885
$ws.Cells.Item(122, 13).Value = -1133.05885

$ws.Cells.Item(131, 8).Value = 1328.7869
$ws.Cells.Item(131, 10).Value = 1143.2632
$ws.Cells.Item(131, 12).Value = 3429.7896
$ws.Cells.Item(131, 14).Value = -13509.7896

$ws.Cells.Item(135, 8).Value = 27783928
$ws.Cells.Item(135, 9).Value = 565.6667
$ws.Cells.Item(135, 10).Value = 55567292
$ws.Cells.Item(135, 11).Value = 5091.0003
$ws.Cells.Item(135, 12).Value = 500105628
$ws.Cells.Item(135, 13).Value = -2556.0003
$ws.Cells.Item(135, 14).Value = -500110698

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 6776.727
$ws.Cells.Item(122, 9).Value = 11060
$ws.Cells.Item(122, 10).Value = 3207.3333
$ws.Cells.Item(122, 11).Value = 33180
$ws.Cells.Item(122, 12).Value = 9621.999899999999
$ws.Cells.Item(122, 13).Value = -30730
$ws.Cells.Item(122, 14).Value = -14521.9999

$ws.Cells.Item(132, 8).Value = 4135.5
$ws.Cells.Item(132, 9).Value = 1589.7354
$ws.Cells.Item(132, 10).Value = 25774.5
$ws.Cells.Item(132, 11).Value = 4769.206200000001
$ws.Cells.Item(132, 12).Value = 77323.5
$ws.Cells.Item(132, 13).Value = -2239.206200000001
$ws.Cells.Item(132, 14).Value = -82383.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3230.2263
$ws.Cells.Item(132, 9).Value = 3222.652
$ws.Cells.Item(132, 10).Value = 3280
$ws.Cells.Item(132, 11).Value = 9667.956
$ws.Cells.Item(132, 12).Value = 9840
$ws.Cells.Item(132, 13).Value = -7137.956
$ws.Cells.Item(132, 14).Value = -14900

$ws.Cells.Item(136, 8).Value = 2931.2307
$ws.Cells.Item(136, 9).Value = 1532.6072
$ws.Cells.Item(136, 10).Value = 6491.364
$ws.Cells.Item(136, 11).Value = 4597.821599999999
$ws.Cells.Item(136, 12).Value = 19474.092
$ws.Cells.Item(136, 13).Value = -2047.821599999999
$ws.Cells.Item(136, 14).Value = -24574.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2253.4827
$ws.Cells.Item(81, 9).Value = 1563.9445
$ws.Cells.Item(81, 10).Value = 3381.818
$ws.Cells.Item(81, 11).Value = 3127.889
$ws.Cells.Item(81, 12).Value = 6763.636
$ws.Cells.Item(81, 13).Value = -2066.889
$ws.Cells.Item(81, 14).Value = -8885.636

$ws.Cells.Item(84, 8).Value = 2253.4827
$ws.Cells.Item(84, 9).Value = 1563.9445
$ws.Cells.Item(84, 10).Value = 3381.818
$ws.Cells.Item(84, 11).Value = 15639.445
$ws.Cells.Item(84, 12).Value = 33818.18
$ws.Cells.Item(84, 13).Value = -10335.445
$ws.Cells.Item(84, 14).Value = -44426.18

$ws.Cells.Item(107, 8).Value = 2612.4707
$ws.Cells.Item(107, 9).Value = 508.15384
$ws.Cells.Item(107, 10).Value = 9451.5
$ws.Cells.Item(107, 11).Value = 1524.46152
$ws.Cells.Item(107, 12).Value = 28354.5
$ws.Cells.Item(107, 13).Value = 395.5384799999999
$ws.Cells.Item(107, 14).Value = -32194.5

$ws.Cells.Item(136, 8).Value = 3952.4353
$ws.Cells.Item(136, 9).Value = 2546.0327
$ws.Cells.Item(136, 10).Value = 7527.0415
$ws.Cells.Item(136, 11).Value = 7638.098100000001
$ws.Cells.Item(136, 12).Value = 22581.1245
$ws.Cells.Item(136, 13).Value = -5088.098100000001
$ws.Cells.Item(136, 14).Value = -27681.1245
